$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Octubre de 2020 a las 14:33"

# --- Reorder Montserrat / Islas Malvinas (rows 215-216) ---
# In the source data these two countries swap places (Montserrat now listed
# before Islas Malvinas). Columns B, C, E, F, G are identical between the two
# rows, so the only visible effect is swapping the country label together
# with the D (Casos activos) and H (Muertes) values.
$ws.Range("A215").Value = "Montserrat"
$ws.Range("A216").Value = "Islas Malvinas"

$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0

# --- Updated country statistics ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 7680028
$ws.Range("C4").Value = 384
$ws.Range("E4").Value = 2569698

# Row 12: Mexico
$ws.Range("B12").Value = 789780
$ws.Range("C12").Value = 3417
$ws.Range("D12").Value = 553937
$ws.Range("E12").Value = 153966
$ws.Range("G12").Value = 180
$ws.Range("H12").Value = 81877

# Row 40: Kuwait
$ws.Range("B40").Value = 108268
$ws.Range("C40").Value = 676
$ws.Range("D40").Value = 100179
$ws.Range("E40").Value = 7457
$ws.Range("G40").Value = 4
$ws.Range("H40").Value = 632

# Row 57: Barein
$ws.Range("E57").Value = 4921
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = 262

# Row 71: Azerbaiyan
$ws.Range("B71").Value = 40931
$ws.Range("C71").Value = 143
$ws.Range("D71").Value = 38713
$ws.Range("E71").Value = 1618
$ws.Range("G71").Value = 2
$ws.Range("H71").Value = 600

# Row 78: Dinamarca
$ws.Range("B78").Value = 30379
$ws.Range("C78").Value = 322
$ws.Range("D78").Value = 23655
$ws.Range("E78").Value = 6061
$ws.Range("G78").Value = 4
$ws.Range("H78").Value = 663

# Row 92: Madagascar
$ws.Range("B92").Value = 16600
$ws.Range("C92").Value = 30
$ws.Range("D92").Value = 15698
$ws.Range("E92").Value = 668
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 234

# Row 148: Islandia
$ws.Range("B148").Value = 3081
$ws.Range("C148").Value = 101
$ws.Range("D148").Value = 2324
$ws.Range("E148").Value = 747
